# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these value cells are stored as text (preserves values such as
# "1.20", "0.970", trailing-zero decimals, and dotted price groupings)
# exactly as scraped, instead of Excel auto-coercing them to numbers.
$textCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'E28', 'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'D32', 'E32', 'E33', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'E38', 'D39', 'E39', 'E40', 'E41', 'D42', 'E42', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'B45', 'C45', 'D45', 'E45', 'E46', 'D47', 'E47', 'E48', 'D49', 'E49', 'E50', 'D51', 'E51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '96.921.05'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.665.74'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D4').Value = '2.68'
$ws.Range('E4').Value = '  +40.72%  '
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '228.19'
$ws.Range('E6').Value = '  -4.44%  '
$ws.Range('D7').Value = '650.04'
$ws.Range('E7').Value = '  -1.77%  '
$ws.Range('D8').Value = '0.429'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '1.20'
$ws.Range('E9').Value = '  +12.38%  '
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.662.45'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '49.45'
$ws.Range('E12').Value = '  +10.25%  '
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '0.0000296'
$ws.Range('E14').Value = '  -8.79%  '
$ws.Range('D15').Value = '6.71'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').Value = '4.352.81'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '96.712.46'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '21.41'
$ws.Range('E18').Value = '  +13.52%  '
$ws.Range('D19').Value = '8.93'
$ws.Range('E19').Value = '  -2.33%  '
$ws.Range('D20').Value = '14.21'
$ws.Range('E20').Value = '  +7.75%  '
$ws.Range('D21').Value = '3.659.99'
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').Value = '0.565'
$ws.Range('E22').Value = '  +11.55%  '
$ws.Range('E23').Value = '  +43.08%  '
$ws.Range('D24').Value = '529.73'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '3.29'
$ws.Range('E25').Value = '  -5.81%  '
$ws.Range('D26').Value = '123.33'
$ws.Range('E26').Value = '  +13.50%  '
$ws.Range('D27').Value = '0.0000204'
$ws.Range('E27').Value = '  -9.84%  '
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').Value = '3.840.91'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Value = '13.08'
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('D31').Value = '13.15'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '3.04'
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('D35').Value = '33.20'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '0.624'
$ws.Range('E36').Value = '  +4.28%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  -4.17%  '
$ws.Range('D39').Value = '606.12'
$ws.Range('E39').Value = '  -7.03%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('D42').Value = '7.12'
$ws.Range('E42').Value = '  +3.43%  '
$ws.Range('D43').Value = '41.80'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0505'
$ws.Range('E44').Value = '  +9.51%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '0.492'
$ws.Range('E45').Value = '  +2.62%  '
$ws.Range('E46').Value = '  -5.29%  '
$ws.Range('D47').Value = '0.970'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('E48').Value = '  -4.23%  '
$ws.Range('D49').Value = '235.17'
$ws.Range('E49').Value = '  +13.28%  '
$ws.Range('E50').Value = '  -3.85%  '
$ws.Range('D51').Value = '8.91'
$ws.Range('E51').Value = '  +1.18%  '
